# Updated cryptos list on Tue Mar  7 17:43:28 UTC 2023 with GitHub Actions
#
# D-column "Price" values are stored as text in the source sheet (they use
# localized dotted-thousands separators and retain trailing/leading zeros
# that a numeric cell would normalize away). Prefixing with an apostrophe
# forces Excel to keep them as literal text, matching the original
# inlineStr storage instead of letting auto-detection coerce them to
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'22.292.16"
$ws.Range("E2").Value = "  -1.16%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.561.34"
$ws.Range("E3").Value = "  -1.03%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5 - USDC
$ws.Range("D5").Value = "'0.9999"
$ws.Range("E5").Value = "  -0.20%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'288.77"
$ws.Range("E6").Value = "  -0.07%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.3776"
$ws.Range("E7").Value = "  +2.21%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3285"
$ws.Range("E8").Value = "  -1.83%  "

# Row 9 - OKB
$ws.Range("D9").Value = "'44.72"
$ws.Range("E9").Value = "  -8.14%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "'1.149"
$ws.Range("E10").Value = "  +0.35%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.07405"
$ws.Range("E11").Value = "  -1.03%  "

# Row 12 - BinanceUSD
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.17%  "

# Row 13 - Solana
$ws.Range("D13").Value = "'20.43"
$ws.Range("E13").Value = "  -2.79%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'5.882"
$ws.Range("E14").Value = "  -2.14%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "'6.782"
$ws.Range("E15").Value = "  -2.67%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "'1.551.85"
$ws.Range("E16").Value = "  -1.62%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "'0.00001080"
$ws.Range("E17").Value = "  -3.54%  "

# Row 18 - TRON
$ws.Range("D18").Value = "'0.06658"
$ws.Range("E18").Value = "  -1.55%  "

# Row 19 - Litecoin (only D changes)
$ws.Range("D19").Value = "'86.44"

# Row 20 - Uniswap
$ws.Range("D20").Value = "'6.427"
$ws.Range("E20").Value = "  +0.01%  "

# Row 21 - Dai (only D changes)
$ws.Range("D21").Value = "'0.9997"

# Row 22 - Avalanche
$ws.Range("D22").Value = "'16.23"
$ws.Range("E22").Value = "  -2.13%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "'11.74"
$ws.Range("E23").Value = "  -3.77%  "

# Row 24 - WrappedBTC
$ws.Range("D24").Value = "'22.279.40"
$ws.Range("E24").Value = "  -1.20%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "'2.297"
$ws.Range("E25").Value = "  -4.32%  "

# Row 26 - LidoDAOToken
$ws.Range("D26").Value = "'2.618"
$ws.Range("E26").Value = "  +0.76%  "

# Row 27 - Monero
$ws.Range("D27").Value = "'151.57"
$ws.Range("E27").Value = "  -0.74%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'19.36"
$ws.Range("E28").Value = "  -1.86%  "

# Row 29 - HuobiToken
$ws.Range("D29").Value = "'4.928"
$ws.Range("E29").Value = "  -1.82%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "'123.21"
$ws.Range("E30").Value = "  -1.11%  "

# Row 31 - WrappedliquidstakedEther2.0
$ws.Range("D31").Value = "'1.727.15"
$ws.Range("E31").Value = "  -1.50%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "'1.088"
$ws.Range("E32").Value = "  +1.68%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "'5.950"
$ws.Range("E33").Value = "  -3.90%  "

# Row 34 - WEMIXTOKEN
$ws.Range("D34").Value = "'1.915"
$ws.Range("E34").Value = "  -4.54%  "

# Row 35 - FraxShare
$ws.Range("D35").Value = "'9.441"
$ws.Range("E35").Value = "  -2.70%  "

# Row 36 - Stellar
$ws.Range("D36").Value = "'0.08231"
$ws.Range("E36").Value = "  -1.23%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "'0.02373"
$ws.Range("E37").Value = "  -3.58%  "

# Row 38 - was InternetComputer(DFINITY), now Hedera (rows 38/39 swapped)
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06329"
$ws.Range("E38").Value = "  -1.06%  "

# Row 39 - was Hedera, now InternetComputer(DFINITY)
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.347"
$ws.Range("E39").Value = "  -1.83%  "

# Row 40 - Algorand
$ws.Range("D40").Value = "'0.2160"
$ws.Range("E40").Value = "  -4.87%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "'1.251"
$ws.Range("E41").Value = "  -3.81%  "

# Row 42 - Aptos
$ws.Range("D42").Value = "'11.09"
$ws.Range("E42").Value = "  -2.72%  "

# Row 43 - TheSandbox
$ws.Range("D43").Value = "'0.6100"
$ws.Range("E43").Value = "  -4.12%  "

# Row 44 - Frax
$ws.Range("D44").Value = "'0.9995"
$ws.Range("E44").Value = "  -0.25%  "

# Row 45 - EnergySwap
$ws.Range("D45").Value = "'13.83"
$ws.Range("E45").Value = "  -1.18%  "

# Row 46 - was Decentraland, now PancakeSwap (rows 46/47 swapped)
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.752"
$ws.Range("E46").Value = "  -0.53%  "

# Row 47 - was PancakeSwap, now Decentraland
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5943"
$ws.Range("E47").Value = "  -4.05%  "

# Row 48 - NEARProtocol
$ws.Range("D48").Value = "'1.995"
$ws.Range("E48").Value = "  -3.40%  "

# Row 49 - Quant
$ws.Range("D49").Value = "'123.12"
$ws.Range("E49").Value = "  -1.49%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "'0.07109"
$ws.Range("E51").Value = "  -2.31%  "
